$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (single '.') -- must be
# protected so Excel keeps storing them as literal text, matching the
# original inlineStr/shared-string cell type.
$textCells = @(
    "D5"
    "D10"
    "D15"
    "D17"
    "D19"
    "D22"
    "D24"
    "D25"
    "D29"
    "D37"
    "D38"
    "D43"
    "D45"
    "D46"
    "D47"
    "D49"
    "D50"
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- value updates ---
$ws.Range("D5").Value = "214.54"
$ws.Range("D10").Value = "19.39"
$ws.Range("D15").Value = "65.09"
$ws.Range("D17").Value = "236.58"
$ws.Range("D19").Value = "7.78"
$ws.Range("D22").Value = "4.40"
$ws.Range("D24").Value = "9.13"
$ws.Range("D25").Value = "145.82"
$ws.Range("D29").Value = "15.71"
$ws.Range("D37").Value = "0.573"
$ws.Range("D38").Value = "0.0168"
$ws.Range("D43").Value = "0.956"
$ws.Range("D45").Value = "0.768"
$ws.Range("D46").Value = "62.20"
$ws.Range("D47").Value = "88.52"
$ws.Range("D49").Value = "0.0504"
$ws.Range("D50").Value = "0.0968"
$ws.Range("D51").Value = "7.48"

$ws.Range("D2").Value = "26.748.91"
$ws.Range("D3").Value = "1.624.86"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").Value = "1.853.53"
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("D13").Value = "1.627.22"
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("E17").Value = "  +11.32%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "26.750.55"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("E19").Value = "  +5.43%  "
$ws.Range("D20").Value = "0.0₃0729"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  +3.26%  "
$ws.Range("E23").Value = "  +3.66%  "
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("E29").Value = "  +3.53%  "
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("D33").Value = "1.476.12"
$ws.Range("E33").Value = "  +10.21%  "
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("E40").Value = "  +3.55%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  +3.16%  "
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("D44").Value = "1.764.37"
$ws.Range("E44").Value = "  +2.40%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("E47").Value = "  +3.24%  "
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("E51").Value = "  +2.01%  "

# Restore default (Normal) style on the protected cells so no stray
# number-format style index is left attached to them.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
